$wb = $excel.ActiveWorkbook

# --- "Formula Test" sheet: add a new row demonstrating bracket characters
#     inside a JETT formula expression ($[...]) ---
$wsFormula = $wb.Worksheets.Item("Formula Test")
$wsFormula.Range("A7").Value = "Bracket Test"
$wsFormula.Range("C7").Value = '$[TEXT(39300.625, "[h]")]'

# --- "Outside Reference" sheet: move the active selection to B1 ---
$wsOutside = $wb.Worksheets.Item("Outside Reference")
$wsOutside.Range("B1").Select()

# Restore "Formula Test" as the active/selected tab (it was active before
# this edit, and the diff does not move the active tab).
$wsFormula.Activate()
$wsFormula.Range("A1:F1").Select()
